# Adapt column header formatting to respective input file names (#7)
#
# 1. Rename the header row (row 1) suffixes:
#      "..._old" -> "..._FV2410"
#      "..._new" -> "..._FV2504"
# 2. Freeze the header row (row 1) via a frozen pane.
# 3. Turn the used range A1:U61 into an Excel Table ("Table1") with an
#    AutoFilter, matching the new header names as its column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
  "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410",
  "diff",
  "Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504",
  "Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Freeze the top (header) row: select the cell below the split, then freeze.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a native Excel table, driving the autofilter +
# table part / relationship plumbing the same way Excel's UI would.
$dataRange = $ws.Range("A1:U61")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)

Write-Host "edit.ps1 applied"
